# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation" on all
#   sheets that show it (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - Shrink the now-narrower "Status" columns to match the tighter autofit
#   width Excel computes for the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
